$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A30").Value = "2024-02-24 23:08:57"
$ws.Range("B30").Value = 21
$ws.Range("C30").Value = 17
$ws.Range("D30").Value = 2
$ws.Range("E30").Value = 7
$ws.Range("F30").Value = 3
$ws.Range("G30").Value = 5
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0.001
$ws.Range("J30").Value = 0.05
$ws.Range("K30").Value = 0.003
$ws.Range("L30").Value = 100
$ws.Range("M30").Value = 500
$ws.Range("N30").Value = 10
$ws.Range("O30").Value = 5
$ws.Range("P30").Value = 3
$ws.Range("Q30").Value = 200
$ws.Range("R30").Value = 6
$ws.Range("S30").Value = 3
$ws.Range("T30").Value = 70
$ws.Range("U30").Value = 0.8095238095238095
$ws.Range("V30").Value = "Data/bombay1.xlsx"
$ws.Range("W30").Value = -18600

$ws.Range("A31").Value = "2024-03-17 21:46:40"
$ws.Range("B31").Value = 54
$ws.Range("C31").Value = 23
$ws.Range("D31").Value = 3
$ws.Range("E31").Value = 7
$ws.Range("F31").Value = 6
$ws.Range("G31").Value = 7
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0.001
$ws.Range("J31").Value = 0.05
$ws.Range("K31").Value = 0.003
$ws.Range("L31").Value = 100
$ws.Range("M31").Value = 500
$ws.Range("N31").Value = 10
$ws.Range("O31").Value = 5
$ws.Range("P31").Value = 3
$ws.Range("Q31").Value = 1000
$ws.Range("R31").Value = 3
$ws.Range("S31").Value = 1
$ws.Range("T31").Value = 70
$ws.Range("U31").Value = 0.4259259259259259
$ws.Range("V31").Value = "Data/bombay1.xlsx"
$ws.Range("W31").Value = -488000

$ws.Range("A32").Value = "2024-03-18 22:01:06"
$ws.Range("B32").Value = 59
$ws.Range("C32").Value = 21
$ws.Range("D32").Value = 7
$ws.Range("E32").Value = 10
$ws.Range("F32").Value = 4
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0.001
$ws.Range("J32").Value = 0.05
$ws.Range("K32").Value = 0.003
$ws.Range("L32").Value = 100
$ws.Range("M32").Value = 500
$ws.Range("N32").Value = 10
$ws.Range("O32").Value = 5
$ws.Range("P32").Value = 2
$ws.Range("Q32").Value = 1000
$ws.Range("R32").Value = 2
$ws.Range("S32").Value = 1
$ws.Range("T32").Value = 70
$ws.Range("U32").Value = 0.3559322033898305
$ws.Range("V32").Value = "Data/bombay1.xlsx"
$ws.Range("W32").Value = 392000

$ws.Range("A33").Value = "2024-03-27 09:03:27"
$ws.Range("B33").Value = 0
$ws.Range("C33").Value = 0
$ws.Range("D33").Value = 0
$ws.Range("E33").Value = 0
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0.001
$ws.Range("J33").Value = 0.05
$ws.Range("K33").Value = 0.003
$ws.Range("L33").Value = 100
$ws.Range("M33").Value = 500
$ws.Range("N33").Value = 10
$ws.Range("O33").Value = 3
$ws.Range("P33").Value = 2
$ws.Range("Q33").Value = 1000
$ws.Range("R33").Value = 6
$ws.Range("S33").Value = 1
$ws.Range("T33").Value = 50
$ws.Range("U33").Value = 0
$ws.Range("V33").Value = "Data/bombay1.xlsx"
$ws.Range("W33").Value = 0

$ws.Range("A34").Value = "2024-03-28 00:08:37"
$ws.Range("B34").Value = 61
$ws.Range("C34").Value = 36
$ws.Range("D34").Value = 14
$ws.Range("E34").Value = 22
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0.001
$ws.Range("J34").Value = 0.05
$ws.Range("K34").Value = 0.003
$ws.Range("L34").Value = 100
$ws.Range("M34").Value = 500
$ws.Range("N34").Value = 10
$ws.Range("O34").Value = 9
$ws.Range("P34").Value = 1
$ws.Range("Q34").Value = 1000
$ws.Range("R34").Value = 7
$ws.Range("S34").Value = 1
$ws.Range("T34").Value = 50
$ws.Range("U34").Value = 0.5901639344262295
$ws.Range("V34").Value = "Data/bombay1.xlsx"
$ws.Range("W34").Value = 1161000
